$wb = $excel.ActiveWorkbook

$wsWater = $wb.Worksheets.Item("Water")
$wsInfra = $wb.Worksheets.Item("Infra")
$wsGlobal = $wb.Worksheets.Item("Global")

# Global sheet: remove the "Water specific cost" row (country-specific water cost
# now lives as "Water demand" on the Water sheet); remaining rows shift up.
$wsGlobal.Rows(2).Delete()

# Water sheet: add new "Water demand (L/kg H2)" parameter row.
$wsWater.Range("A6").Value = "Water demand  (L/kg H2)"
$wsWater.Range("B6").Value = 21
$wsWater.Range("A6:B6").Font.Color = 0

# Update selections / active cells to match the saved view state.
$wsInfra.Range("D31").Select()
$wsGlobal.Range("D11").Select()

$wsWater.Activate()
$wsWater.Range("A6").Select()
